# Applies the LOQ4246.xlsx content re-shuffle described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: B/C text changes (cells already exist with correct styles) ---
$ws.Range("B10").Value = "5840535 - Messias Borges Silva"
$ws.Range("C10").Value = "5840535 - Messias Borges Silva"

# --- Row 13: becomes a full A/B/C row (A13 did not exist before) ---
# Give A13 the label style (s=1) by copying format from an existing column-A label cell.
$ws.Range("A9").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# --- Row 14: text changes only (styles already correct) ---
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "Sustainability. Environment acts and protocols. Environmental issues. Natural resources and their pollution, Carbon credits, Zero waste concept ISO 14000, Life Cycle Analysis, Environmental Impact Assessment studies, Sustainable habitat, Conventional and renewable sources, Technology and sustainable development, Sustainable urbanization, Industrial Ecology."
$ws.Range("C14").Value = "Sustainability. Environment acts and protocols. Environmental issues. Natural resources and their pollution, Carbon credits, Zero waste concept ISO 14000, Life Cycle Analysis, Environmental Impact Assessment studies, Sustainable habitat, Conventional and renewable sources, Technology and sustainable development, Sustainable urbanization, Industrial Ecology."

# --- Row 15: text changes, height 60 -> 120 ---
$ws.Range("A15").Value = "Programa:"
# B15/C15 become the text "01/01/2018"; typing that literal would be auto-parsed as a date by
# Excel, so instead copy the *value* (not format) from B8/C8, which already hold it as text.
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial(-4163)

# --- Row 16: text changes only ---
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "Sustainability- need and concept, challenges,Environment acts and protocols, Global, Regional and Local environmental issues, Natural resources and their pollution, Carbon credits, Zero waste concept  ISO 14000, Life Cycle Analysis, Environmental Impact Assessment studies, Sustainable habitat, Green buildings, Green materials, Energy, Conventional and renewable sources,Technology and sustainable development,Sustainable urbanization, Industrial Ecology."
$ws.Range("C16").Value = "Sustainability- need and concept, challenges,Environment acts and protocols, Global, Regional and Local environmental issues, Natural resources and their pollution, Carbon credits, Zero waste concept  ISO 14000, Life Cycle Analysis, Environmental Impact Assessment studies, Sustainable habitat, Green buildings, Green materials, Energy, Conventional and renewable sources,Technology and sustainable development,Sustainable urbanization, Industrial Ecology."

# --- Row 17: becomes an A-only row (was a full A/B/C row) ---
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").Clear()
$ws.Range("C17").Clear()

# --- Row 18: becomes a full A/B/C row (B18/C18 did not exist before) ---
$ws.Range("A18").Value = "Método:"
# B18/C18 need the value styles (s=2 / s=3); copy formats from an existing B/C pair.
$ws.Range("B14").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("B18").Value = "5840535 - Messias Borges Silva"
$ws.Range("C18").Value = "5840535 - Messias Borges Silva"

# --- Row 19: label text changes only ---
$ws.Range("A19").Value = "Critério:"

# --- Row 20: label text changes only ---
$ws.Range("A20").Value = "Norma de recuperação:"

# --- Row 21: label text changes, height 60 -> 120 ---
$ws.Range("A21").Value = "Bibliografia:"

# --- Row 22: removed entirely ---
$ws.Rows.Item(22).Delete()

# --- Row heights ---
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(17).AutoFit()

Write-Host "Edit applied"
